$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B39: it was stored as a text value "4"; it should become a real number 4
$ws.Range("B39").Value = 4

# Append new row 40 with the annotation data
$ws.Range("A40").Value = "Sunsi Wu"

# B40 keeps the politeness score "3" stored as text (matches source data),
# so force a temporary text format while assigning, then drop the format
# again so the cell is left without any explicit style.
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "3"
$ws.Range("B40").ClearFormats()

$ws.Range("C40").Value = "necessary;sufficient"
$ws.Range("D40").Value = "APC"
$ws.Range("E40").Value = "RES"
$ws.Range("F40").Value = "df7b0ece-3727-4ec6-95ce-2a2839e398ed"
$ws.Range("G40").Value = "SkhQHMW0W_annotated.xlsx"
$ws.Range("H40").Value = "This is necessary to get an idea of total amount of communication that was sufficient to reach perplexity 72.24 at the end of 40-th epoch."
